$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# ---------------------------------------------------------------------
# 2. Split the mailing-address paragraph (the first occurrence only, the
#    one right below the addressee name, not the "PROPERTY ADDRESS" one
#    further down in the table) into two paragraphs:
#      "2970 Lamory Pl"
#      "Santa Clara, CA 95051"
# ---------------------------------------------------------------------
$addrRange = $d.Content
$addrRange.Find.Execute("2970 Lamory Pl, Santa Clara CA 95051", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

$full = $addrRange.Text
$plEndOffset = $full.IndexOf("Pl") + 2
$splitAt = $addrRange.Start + $plEndOffset
$restLen = $full.Length - $plEndOffset

# Insert a paragraph break right after "Pl" (this leaves ", Santa Clara CA 95051"
# as the start of the following paragraph).
$splitPoint = $d.Range($splitAt, $splitAt)
$splitPoint.InsertParagraphAfter()

# The new paragraph mark occupies $splitAt, so the remaining text now begins
# at $splitAt + 1; re-grab that exact bounded range (not a collapsed one) and
# rewrite it so it reads "Santa Clara, CA 95051".
$secondStart = $splitAt + 1
$secondPara = $d.Range($secondStart, $secondStart + $restLen)
$secondPara.Text = "Santa Clara, CA 95051"

# ---------------------------------------------------------------------
# 3. Remove the empty "No Spacing" paragraph that used to sit right after
#    the "Board of Directors" signature line.
# ---------------------------------------------------------------------
$bodRange = $d.Content
$bodRange.Find.Execute("Board of Directors", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0)
$bodParagraph = $bodRange.Paragraphs(1)
$trailingParagraph = $bodParagraph.Next()
if ($trailingParagraph.Range.Text.Trim() -eq "") {
    $trailingParagraph.Range.Delete()
}
